$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 106, shifting existing rows 106-157 down to 107-158
$ws.Rows(106).Insert()

# Populate the newly inserted row 106 with the new data point
$ws.Range("A106").Value2 = 8
$ws.Range("B106").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C106").Value2 = "Coquimbo"
$ws.Range("D106").Value2 = 44489
$ws.Range("E106").Value2 = 4
$ws.Range("F106").Value2 = 100112003
$ws.Range("G106").Value2 = "Ajo"
$ws.Range("H106").Value2 = "Chino"
$ws.Range("I106").Value2 = "Primera"
$ws.Range("J106").Value2 = 640
$ws.Range("K106").Value2 = 17000
$ws.Range("L106").Value2 = 18000
$ws.Range("M106").Value2 = 17500
$ws.Range("N106").Value2 = "$/caja 10 kilos"
$ws.Range("O106").Value2 = "China"
$ws.Range("P106").Value2 = 1750
$ws.Range("Q106").Value2 = 10
$ws.Range("R106").Value2 = "Hortaliza"
